$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.43918213709012
$ws.Range("C2").Value = 8.299527167885106
$ws.Range("D2").Value = 13.516024046038
$ws.Range("E2").Value = 13.9219282601955
$ws.Range("G2").Value = 35.13222257082307
$ws.Range("H2").Value = 15.88121794057369
$ws.Range("I2").Value = 24.34324190708805
$ws.Range("J2").Value = 8.430017636232815
$ws.Range("L2").Value = 12.33806177167017
$ws.Range("N2").Value = 17.93985269201614
$ws.Range("O2").Value = 25.02924558630034

$ws.Range("B3").Value = 17.0062613004529
$ws.Range("C3").Value = 8.059732190455096
$ws.Range("D3").Value = 13.52097337116417
$ws.Range("E3").Value = 13.95636325109047
$ws.Range("G3").Value = 35.15611884589168
$ws.Range("H3").Value = 15.92548084319513
$ws.Range("I3").Value = 24.44330165851723
$ws.Range("J3").Value = 8.439498305988668
$ws.Range("L3").Value = 12.32261436424844
$ws.Range("N3").Value = 17.9832527258697
$ws.Range("O3").Value = 25.09050384396986

$ws.Range("B4").Value = 16.73678669722084
$ws.Range("C4").Value = 7.907693877438902
$ws.Range("D4").Value = 13.52647241683673
$ws.Range("E4").Value = 13.97931693002928
$ws.Range("G4").Value = 35.18140575154056
$ws.Range("H4").Value = 15.95530723337479
$ws.Range("I4").Value = 24.50941584172255
$ws.Range("J4").Value = 8.445632347984388
$ws.Range("L4").Value = 12.3147992175039
$ws.Range("N4").Value = 18.01165699815161
$ws.Range("O4").Value = 25.13354347283185

$ws.Range("B5").Value = 16.62620903013301
$ws.Range("C5").Value = 7.84459105393269
$ws.Range("D5").Value = 13.52933270341052
$ws.Range("E5").Value = 13.98912632399854
$ws.Range("G5").Value = 35.19437445016268
$ws.Range("H5").Value = 15.96812739762032
$ws.Range("I5").Value = 24.53753297611201
$ws.Range("J5").Value = 8.448210918640973
$ws.Range("L5").Value = 12.3120370256487
$ws.Range("N5").Value = 18.02367445604103
$ws.Range("O5").Value = 25.1524444301878

$ws.Range("B6").Value = 16.60780621131816
$ws.Range("C6").Value = 7.834045572918101
$ws.Range("D6").Value = 13.52984507848806
$ws.Range("E6").Value = 13.99078269208896
$ws.Range("G6").Value = 35.19668861644297
$ws.Range("H6").Value = 15.97029636889825
$ws.Range("I6").Value = 24.54227274654912
$ws.Range("J6").Value = 8.448643860229552
$ws.Range("L6").Value = 12.31160396106257
$ws.Range("N6").Value = 18.0256966947998
$ws.Range("O6").Value = 25.15566509659777

$ws.Range("B7").Value = 16.73529829264894
$ws.Range("C7").Value = 7.906847405194366
$ws.Range("D7").Value = 13.52650848297711
$ws.Range("E7").Value = 13.97944737786862
$ws.Range("G7").Value = 35.1815698731142
$ws.Range("H7").Value = 15.95547743587388
$ws.Range("I7").Value = 24.50979028277431
$ws.Range("J7").Value = 8.445666803732777
$ws.Range("L7").Value = 12.3147602513459
$ws.Range("N7").Value = 18.01181727710581
$ws.Range("O7").Value = 25.13379286692827

$ws.Range("B8").Value = 17.29076149045172
$ws.Range("C8").Value = 8.217878259817226
$ws.Range("D8").Value = 13.51722049252617
$ws.Range("E8").Value = 13.93342588561553
$ws.Range("G8").Value = 35.13825645382479
$ws.Range("H8").Value = 15.89592980618027
$ws.Range("I8").Value = 24.37677116198918
$ws.Range("J8").Value = 8.433221769862035
$ws.Range("L8").Value = 12.33239065933975
$ws.Range("N8").Value = 17.95445294938264
$ws.Range("O8").Value = 25.04923956795761

$ws.Range("B9").Value = 18.34456212777544
$ws.Range("C9").Value = 8.787244054413776
$ws.Range("D9").Value = 13.51848521785585
$ws.Range("E9").Value = 13.85752827226879
$ws.Range("G9").Value = 35.13768966434704
$ws.Range("H9").Value = 15.80019091138756
$ws.Range("I9").Value = 24.15307753462168
$ws.Range("J9").Value = 8.411289332309501
$ws.Range("L9").Value = 12.38008485658335
$ws.Range("N9").Value = 17.8558626521284
$ws.Range("O9").Value = 24.92660060891965

$ws.Range("B10").Value = 19.08919173436037
$ws.Range("C10").Value = 9.177899749102011
$ws.Range("D10").Value = 13.53122066685014
$ws.Range("E10").Value = 13.81049440712235
$ws.Range("G10").Value = 35.18880115032898
$ws.Range("H10").Value = 15.74269640669241
$ws.Range("I10").Value = 24.01144422569995
$ws.Range("J10").Value = 8.39666830518475
$ws.Range("L10").Value = 12.4229433343857
$ws.Range("N10").Value = 17.79185385198298
$ws.Range("O10").Value = 24.86295371706353

$ws.Range("B11").Value = 19.42001515733451
$ws.Range("C11").Value = 9.349090863513377
$ws.Range("D11").Value = 13.5395583494702
$ws.Range("E11").Value = 13.79098830754273
$ws.Range("G11").Value = 35.2232230194069
$ws.Range("H11").Value = 15.71933399540035
$ws.Range("I11").Value = 23.95195779908145
$ws.Range("J11").Value = 8.390337891590818
$ws.Range("L11").Value = 12.44409516660818
$ws.Range("N11").Value = 17.76455401647082
$ws.Range("O11").Value = 24.83976808720812

$ws.Range("B12").Value = 19.54403940076845
$ws.Range("C12").Value = 9.412941180453583
$ws.Range("D12").Value = 13.54307935145876
$ws.Range("E12").Value = 13.78387326383896
$ws.Range("G12").Value = 35.23785919568579
$ws.Range("H12").Value = 15.71088906918066
$ws.Range("I12").Value = 23.93014397480624
$ws.Range("J12").Value = 8.387986623906233
$ws.Range("L12").Value = 12.45233880710873
$ws.Range("N12").Value = 17.75447697920378
$ws.Range("O12").Value = 24.83181917043007

$ws.Range("B13").Value = 19.51738592552551
$ws.Range("C13").Value = 9.399233850283919
$ws.Range("D13").Value = 13.54230490208552
$ws.Range("E13").Value = 13.7853935453906
$ws.Range("G13").Value = 35.23463591304834
$ws.Range("H13").Value = 15.7126899539347
$ws.Range("I13").Value = 23.93481026322359
$ws.Range("J13").Value = 8.388490972268361
$ws.Range("L13").Value = 12.45055305682822
$ws.Range("N13").Value = 17.75663565922988
$ws.Range("O13").Value = 24.83349413180888

$ws.Range("B14").Value = 19.43024440352283
$ws.Range("C14").Value = 9.354363634373728
$ws.Range("D14").Value = 13.53984074858127
$ws.Range("E14").Value = 13.79039750903796
$ws.Range("G14").Value = 35.22439508848834
$ws.Range("H14").Value = 15.71863116913221
$ws.Range("I14").Value = 23.95014887619543
$ws.Range("J14").Value = 8.390143532183531
$ws.Range("L14").Value = 12.44476871255766
$ws.Range("N14").Value = 17.76371974886355
$ws.Range("O14").Value = 24.83909746296804

$ws.Range("B15").Value = 19.37670152547264
$ws.Range("C15").Value = 9.326751151711932
$ws.Range("D15").Value = 13.53837867965858
$ws.Range("E15").Value = 13.79349792965325
$ws.Range("G15").Value = 35.21833064986325
$ws.Range("H15").Value = 15.7223226871994
$ws.Range("I15").Value = 23.95963704087206
$ws.Range("J15").Value = 8.391161747324716
$ws.Range("L15").Value = 12.44125596928692
$ws.Range("N15").Value = 17.76809290638423
$ws.Range("O15").Value = 24.84263792502926

$ws.Range("B16").Value = 19.06740249838616
$ws.Range("C16").Value = 9.166577422304997
$ws.Range("D16").Value = 13.53072679815443
$ws.Range("E16").Value = 13.81180718268266
$ws.Range("G16").Value = 35.18677600559016
$ws.Range("H16").Value = 15.74427941567403
$ws.Range("I16").Value = 24.01543140761195
$ws.Range("J16").Value = 8.397088449813729
$ws.Range("L16").Value = 12.42159397311958
$ws.Range("N16").Value = 17.79367449153777
$ws.Range("O16").Value = 24.86458515600805

$ws.Range("B17").Value = 18.87554714519948
$ws.Range("C17").Value = 9.066617362460361
$ws.Range("D17").Value = 13.52668281032999
$ws.Range("E17").Value = 13.82352316510202
$ws.Range("G17").Value = 35.17027647139782
$ws.Range("H17").Value = 15.75846459045769
$ws.Range("I17").Value = 24.05092663784375
$ws.Range("J17").Value = 8.400806300382559
$ws.Range("L17").Value = 12.40995306667013
$ws.Range("N17").Value = 17.80983317406545
$ws.Range("O17").Value = 24.87952754035705

$ws.Range("B18").Value = 18.7644584562982
$ws.Range("C18").Value = 9.008511297391941
$ws.Range("D18").Value = 13.52459643944246
$ws.Range("E18").Value = 13.83043978549912
$ws.Range("G18").Value = 35.16183826795537
$ws.Range("H18").Value = 15.7668863140679
$ws.Range("I18").Value = 24.07180774553239
$ws.Range("J18").Value = 8.402974913641833
$ws.Range("L18").Value = 12.40341359023915
$ws.Range("N18").Value = 17.81929838395805
$ws.Range("O18").Value = 24.88866480789916

$ws.Range("B19").Value = 18.72672250479602
$ws.Range("C19").Value = 8.988733765972825
$ws.Range("D19").Value = 13.52393124725011
$ws.Range("E19").Value = 13.83281219623237
$ws.Range("G19").Value = 35.1591620213242
$ws.Range("H19").Value = 15.76978288573125
$ws.Range("I19").Value = 24.07895757954379
$ws.Range("J19").Value = 8.403714362602637
$ws.Range("L19").Value = 12.40122636250706
$ws.Range("N19").Value = 17.82253255804085
$ws.Range("O19").Value = 24.89185170304251

$ws.Range("B20").Value = 18.89604776702906
$ws.Range("C20").Value = 9.077321880173566
$ws.Range("D20").Value = 13.52708851437417
$ws.Range("E20").Value = 13.82225756896885
$ws.Range("G20").Value = 35.17192403813013
$ws.Range("H20").Value = 15.75692735390847
$ws.Range("I20").Value = 24.04709995333011
$ws.Range("J20").Value = 8.400407404324566
$ws.Range("L20").Value = 12.41117613811828
$ws.Range("N20").Value = 17.80809534453762
$ws.Range("O20").Value = 24.87788070737928

$ws.Range("B21").Value = 19.45587480201458
$ws.Range("C21").Value = 9.367569879856658
$ws.Range("D21").Value = 13.54055467792331
$ws.Range("E21").Value = 13.78892035752841
$ws.Range("G21").Value = 35.22735965613148
$ws.Range("H21").Value = 15.71687517938722
$ws.Range("I21").Value = 23.94562420706316
$ws.Range("J21").Value = 8.389656890562184
$ws.Range("L21").Value = 12.44646140064235
$ws.Range("N21").Value = 17.76163190635326
$ws.Range("O21").Value = 24.83742906509664

$ws.Range("B22").Value = 19.81441393856004
$ws.Range("C22").Value = 9.551559637270184
$ws.Range("D22").Value = 13.55147434467286
$ws.Range("E22").Value = 13.76871490184483
$ws.Range("G22").Value = 35.27292028287319
$ws.Range("H22").Value = 15.69304157883457
$ws.Range("I22").Value = 23.88345705688264
$ws.Range("J22").Value = 8.382898387189055
$ws.Range("L22").Value = 12.47088364843852
$ws.Range("N22").Value = 17.73278533465422
$ws.Range("O22").Value = 24.81583567359583

$ws.Range("B23").Value = 19.62376157004511
$ws.Range("C23").Value = 9.453894409558803
$ws.Range("D23").Value = 13.54545322818591
$ws.Range("E23").Value = 13.77935423928216
$ws.Range("G23").Value = 35.24775216581219
$ws.Range("H23").Value = 15.70554752838367
$ws.Range("I23").Value = 23.91625630504734
$ws.Range("J23").Value = 8.38648111077997
$ws.Range("L23").Value = 12.45772589554888
$ws.Range("N23").Value = 17.74804241596837
$ws.Range("O23").Value = 24.82691679152941

$ws.Range("B24").Value = 18.88678188311691
$ws.Range("C24").Value = 9.072484351065679
$ws.Range("D24").Value = 13.52690435229995
$ws.Range("E24").Value = 13.82282918133563
$ws.Range("G24").Value = 35.17117590964721
$ws.Range("H24").Value = 15.75762150850971
$ws.Range("I24").Value = 24.04882851968942
$ws.Range("J24").Value = 8.400587648120892
$ws.Range("L24").Value = 12.41062271051809
$ws.Range("N24").Value = 17.80888047095193
$ws.Range("O24").Value = 24.87862353776879

$ws.Range("B25").Value = 18.06413741851732
$ws.Range("C25").Value = 8.637879596343923
$ws.Range("D25").Value = 13.51606382586664
$ws.Range("E25").Value = 13.8765263476789
$ws.Range("G25").Value = 35.12879351601034
$ws.Range("H25").Value = 15.82383659807592
$ws.Range("I25").Value = 24.2096083358976
$ws.Range("J25").Value = 8.416959463711175
$ws.Range("L25").Value = 12.44125596928692
$ws.Range("N25").Value = 17.76809290638423
$ws.Range("O25").Value = 24.84263792502926

